$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A31").Value = 112017488
$ws.Range("Q31").Value = 682955.8308828628
$ws.Range("R31").Value = 6575473.896637772
$ws.Range("A32").Value = 112017430
$ws.Range("B32").Value = 90709
$ws.Range("D32").Value = "NT"
$ws.Range("E32").Value = 5448
$ws.Range("F32").Value = "Svartvit taggsvamp"
$ws.Range("G32").Value = "Phellodon connatus"
$ws.Range("H32").Value = "(Schultz) nom.prov"
$ws.Range("Q32").Value = 682793.1335561723
$ws.Range("R32").Value = 6575519.79500053
$ws.Range("A33").Value = 112017534
$ws.Range("B33").Value = 87992
$ws.Range("D33").Value = "VU"
$ws.Range("E33").Value = 1593
$ws.Range("F33").Value = "Lakritsmusseron"
$ws.Range("G33").Value = "Tricholoma apium"
$ws.Range("H33").Value = "Jul.Schäff."
$ws.Range("I33").Value = "4"
$ws.Range("J33").Value = "fruktkroppar"
$ws.Range("Q33").Value = 683072.5368938858
$ws.Range("R33").Value = 6575477.991881827
$ws.Range("A34").Value = 112017512
$ws.Range("B34").Value = 88032
$ws.Range("D34").Value = "VU"
$ws.Range("E34").Value = 6276
$ws.Range("F34").Value = "Goliatmusseron"
$ws.Range("G34").Value = "Tricholoma matsutake"
$ws.Range("H34").Value = "(S.Ito & S.Imai) Singer"
$ws.Range("I34").Value = "4"
$ws.Range("J34").Value = "fruktkroppar"
$ws.Range("Q34").Value = 683036.8460961942
$ws.Range("R34").Value = 6575484.458868909
$ws.Range("A35").Value = 112017326
$ws.Range("B35").Value = 90660
$ws.Range("D35").Value = "NT"
$ws.Range("E35").Value = 4362
$ws.Range("F35").Value = "Blå taggsvamp"
$ws.Range("G35").Value = "Hydnellum caeruleum"
$ws.Range("H35").Value = "(Hornem.) P.Karst."
$ws.Range("Q35").Value = 682713.7813606198
$ws.Range("R35").Value = 6575496.010644327
$ws.Range("A36").Value = 112017465
$ws.Range("I36").Value = "3"
$ws.Range("Q36").Value = 682896.4696766059
$ws.Range("R36").Value = 6575514.027787391
$ws.Range("A37").Value = 112017392
$ws.Range("B37").Value = 90710
$ws.Range("D37").Value = "NT"
$ws.Range("E37").Value = 5449
$ws.Range("F37").Value = "Svart taggsvamp"
$ws.Range("G37").Value = "Phellodon niger"
$ws.Range("H37").Value = "(Fr.:Fr.) P.Karst."
$ws.Range("Q37").Value = 682712.0453105029
$ws.Range("R37").Value = 6575457.539765021
$ws.Range("A38").Value = 112017413
$ws.Range("B38").Value = 90709
$ws.Range("E38").Value = 5448
$ws.Range("F38").Value = "Svartvit taggsvamp"
$ws.Range("G38").Value = "Phellodon connatus"
$ws.Range("H38").Value = "(Schultz) nom.prov"
$ws.Range("Q38").Value = 682733.9332997696
$ws.Range("R38").Value = 6575482.138353716
$ws.Range("A39").Value = 112017130
$ws.Range("B39").Value = 90666
$ws.Range("D39").Value = "LC"
$ws.Range("E39").Value = 4364
$ws.Range("F39").Value = "Dropptaggsvamp"
$ws.Range("G39").Value = "Hydnellum ferrugineum"
$ws.Range("H39").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("Q39").Value = 682695.3118543178
$ws.Range("R39").Value = 6575453.662799283
$ws.Range("A40").Value = 112017447
$ws.Range("B40").Value = 90666
$ws.Range("D40").Value = "LC"
$ws.Range("E40").Value = 4364
$ws.Range("F40").Value = "Dropptaggsvamp"
$ws.Range("G40").Value = "Hydnellum ferrugineum"
$ws.Range("H40").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("Q40").Value = 682844.1942409466
$ws.Range("R40").Value = 6575513.554896963
$ws.Range("A41").Value = 112017252
$ws.Range("B41").Value = 90666
$ws.Range("E41").Value = 4364
$ws.Range("F41").Value = "Dropptaggsvamp"
$ws.Range("G41").Value = "Hydnellum ferrugineum"
$ws.Range("H41").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("Q41").Value = 682710.810501094
$ws.Range("R41").Value = 6575493.820233095
$ws.Range("A42").Value = 112017159
$ws.Range("B42").Value = 90710
$ws.Range("D42").Value = "NT"
$ws.Range("E42").Value = 5449
$ws.Range("F42").Value = "Svart taggsvamp"
$ws.Range("G42").Value = "Phellodon niger"
$ws.Range("H42").Value = "(Fr.:Fr.) P.Karst."
$ws.Range("Q42").Value = 682698.5384611045
$ws.Range("R42").Value = 6575482.480741166
$ws.Range("A43").Value = 112017224
$ws.Range("B43").Value = 90678
$ws.Range("E43").Value = 4366
$ws.Range("F43").Value = "Skarp dropptaggsvamp"
$ws.Range("G43").Value = "Hydnellum peckii"
$ws.Range("Q43").Value = 682702.748818734
$ws.Range("R43").Value = 6575490.872789856

# Clear cells that become empty
$ws.Range("I37").Value = ""
$ws.Range("J37").Value = ""
$ws.Range("I42").Value = ""
$ws.Range("J42").Value = ""
